# "Generate Report for Handback"
#
# The localization status workbook gets refreshed after a handback run:
#   - Status goes from "Ready for handoff" to "Handed back: in sync with en-US"
#     (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3 — they all share one string).
#   - The per-language tables (zh-cn / de-de) get their "Latest Target File"
#     (col I) and "Latest Handback File" (col J) populated, plus a real
#     "Latest Handback DateTime" (col K) instead of the 0001-01-01 placeholder.
#   - Column I/J widen out to fit the long file names, and the Status columns
#     widen a bit for the longer text.

$wb = $excel.ActiveWorkbook

$missing = [System.Reflection.Missing]::Value

$newStatus = "Handed back: in sync with en-US"

$zhcnTargetMd96 = "96c4c545-2127-42e8-a97e-09db01a99ce8.md"
$zhcnTargetMd28 = "f28701aa-d03d-4e95-961b-6f4317b50766.md"
$url96 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6d40bc3c8c14a832894e2697c4b9cb0f1fb239e2/e2e/96c4c545-2127-42e8-a97e-09db01a99ce8.md"
$url28 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6d40bc3c8c14a832894e2697c4b9cb0f1fb239e2/e2e/f28701aa-d03d-4e95-961b-6f4317b50766.md"

# ---------------------------------------------------------------------------
# Overview sheet: status text + wider status columns
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsOverview.Columns.Item(5).ColumnWidth = 29.1
$wsOverview.Columns.Item(6).ColumnWidth = 29.1

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $url96, $missing, $missing, $zhcnTargetMd96) | Out-Null
$wsZh.Range("J2").Value = "96c4c545-2127-42e8-a97e-09db01a99ce8.684e6766958083b3a9c072393596464a9b38fa44.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-05 00:54:34"

$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $url28, $missing, $missing, $zhcnTargetMd28) | Out-Null
$wsZh.Range("J3").Value = "f28701aa-d03d-4e95-961b-6f4317b50766.5cf3bcd98dd282b9cb3e319ff4a70c65ada9717e.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-09-05 00:54:34"

$wsZh.Columns.Item(3).ColumnWidth = 29.1
$wsZh.Columns.Item(9).ColumnWidth = 39.17
$wsZh.Columns.Item(10).ColumnWidth = 39.17

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $url96, $missing, $missing, $zhcnTargetMd96) | Out-Null
$wsDe.Range("J2").Value = "96c4c545-2127-42e8-a97e-09db01a99ce8.684e6766958083b3a9c072393596464a9b38fa44.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-05 00:54:41"

$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $url28, $missing, $missing, $zhcnTargetMd28) | Out-Null
$wsDe.Range("J3").Value = "f28701aa-d03d-4e95-961b-6f4317b50766.5cf3bcd98dd282b9cb3e319ff4a70c65ada9717e.de-de.xlf"
$wsDe.Range("K3").Value = "2016-09-05 00:54:41"

$wsDe.Columns.Item(3).ColumnWidth = 29.1
$wsDe.Columns.Item(9).ColumnWidth = 39.17
$wsDe.Columns.Item(10).ColumnWidth = 39.17
